$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 8's formatting (fills/borders/fonts/number formats) down into row 9,
# mirroring the "copy the last entry to start a new one" edit.
$ws.Range("A8:N8").Copy($ws.Range("A9:N9"))

# New meeting entry values for row 9 ("Offline 7")
$ws.Range("A9").Value = "Offline 7"
$ws.Range("B9").Value = 6
$ws.Range("C9").Value = "17/10/2020"
$ws.Range("N9").Value = 5

# Đỗ Trường Giang (column E) is available for this meeting -> unmark it
$ws.Range("E9").Interior.Color = $ws.Range("L8").Interior.Color

# Restore the selection to reflect where editing finished
$ws.Range("B10").Select()
